# Updates league database rows: swaps/rotates the match-data columns (B:AD)
# between certain rows while keeping the row index in column A fixed.
# This mirrors a re-ordering of underlying source records for:
#   - rows 87 and 88            (full swap)
#   - rows 119, 120 and 121     (3-way rotation: 119<-120, 120<-121, 121<-119)
#   - rows 226 and 227          (full swap)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowData($row) {
    return [ordered]@{
        B  = $ws.Cells.Item($row, 2).Value2
        C  = $ws.Cells.Item($row, 3).Value2
        D  = $ws.Cells.Item($row, 4).Value2
        E  = $ws.Cells.Item($row, 5).Value2
        F  = $ws.Cells.Item($row, 6).Value2
        G  = $ws.Cells.Item($row, 7).Value2
        H  = $ws.Cells.Item($row, 8).Value2
        I  = $ws.Cells.Item($row, 9).Value2
        J  = $ws.Cells.Item($row, 10).Value2
        K  = $ws.Cells.Item($row, 11).Value2
        L  = $ws.Cells.Item($row, 12).Value2
        M  = $ws.Cells.Item($row, 13).Value2
        N  = $ws.Cells.Item($row, 14).Value2
        O  = $ws.Cells.Item($row, 15).Value2
        P  = $ws.Cells.Item($row, 16).Value2
        Q  = $ws.Cells.Item($row, 17).Value2
        R  = $ws.Cells.Item($row, 18).Value2
        S  = $ws.Cells.Item($row, 19).Value2
        T  = $ws.Cells.Item($row, 20).Value2
        U  = $ws.Cells.Item($row, 21).Value2
        V  = $ws.Cells.Item($row, 22).Value2
        W  = $ws.Cells.Item($row, 23).Value2
        X  = $ws.Cells.Item($row, 24).Value2
        Y  = $ws.Cells.Item($row, 25).Value2
        Z  = $ws.Cells.Item($row, 26).Value2
        AA = $ws.Cells.Item($row, 27).Value2
        AB = $ws.Cells.Item($row, 28).Value2
        AC = $ws.Cells.Item($row, 29).Value2
        AD = $ws.Cells.Item($row, 30).Value2
    }
}

function Set-RowData($row, $data) {
    $ws.Cells.Item($row, 2).Value  = $data.B
    $ws.Cells.Item($row, 3).Value  = $data.C
    $ws.Cells.Item($row, 4).Value  = $data.D
    $ws.Cells.Item($row, 5).Value  = $data.E
    $ws.Cells.Item($row, 6).Value  = $data.F
    $ws.Cells.Item($row, 7).Value  = $data.G
    $ws.Cells.Item($row, 8).Value  = $data.H
    $ws.Cells.Item($row, 9).Value  = $data.I
    $ws.Cells.Item($row, 10).Value = $data.J
    $ws.Cells.Item($row, 11).Value = $data.K
    $ws.Cells.Item($row, 12).Value = $data.L
    $ws.Cells.Item($row, 13).Value = $data.M
    $ws.Cells.Item($row, 14).Value = $data.N
    $ws.Cells.Item($row, 15).Value = $data.O
    $ws.Cells.Item($row, 16).Value = $data.P
    $ws.Cells.Item($row, 17).Value = $data.Q
    $ws.Cells.Item($row, 18).Value = $data.R
    $ws.Cells.Item($row, 19).Value = $data.S
    $ws.Cells.Item($row, 20).Value = $data.T
    $ws.Cells.Item($row, 21).Value = $data.U
    $ws.Cells.Item($row, 22).Value = $data.V
    $ws.Cells.Item($row, 23).Value = $data.W
    $ws.Cells.Item($row, 24).Value = $data.X
    $ws.Cells.Item($row, 25).Value = $data.Y
    $ws.Cells.Item($row, 26).Value = $data.Z
    $ws.Cells.Item($row, 27).Value = $data.AA
    $ws.Cells.Item($row, 28).Value = $data.AB
    $ws.Cells.Item($row, 29).Value = $data.AC
    $ws.Cells.Item($row, 30).Value = $data.AD
}

# --- Pair swap: rows 87 and 88 ---
$row87 = Get-RowData 87
$row88 = Get-RowData 88
Set-RowData 87 $row88
Set-RowData 88 $row87

# --- 3-way rotation: rows 119, 120, 121 ---
# new119 = old120 ; new120 = old121 ; new121 = old119
$row119 = Get-RowData 119
$row120 = Get-RowData 120
$row121 = Get-RowData 121
Set-RowData 119 $row120
Set-RowData 120 $row121
Set-RowData 121 $row119

# --- Pair swap: rows 226 and 227 ---
$row226 = Get-RowData 226
$row227 = Get-RowData 227
Set-RowData 226 $row227
Set-RowData 227 $row226
